# Auto-generated Excel COM-interop script applying the scheduled-runner price/profit updates
# to the FFXIV Leve profit tracking sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 214.81818
$ws.Range("I4").Value = 214.81818
$ws.Range("K4").Value = 214.81818
$ws.Range("M4").Value = -100.81818

$ws.Range("H62").Value = 2349.4
$ws.Range("I62").Value = 2842
$ws.Range("J62").Value = 1200
$ws.Range("K62").Value = 2842
$ws.Range("L62").Value = 1200
$ws.Range("M62").Value = -2218
$ws.Range("N62").Value = -2448

$ws.Range("H64").Value = 3478.25
$ws.Range("I64").Value = 3471
$ws.Range("K64").Value = 3471
$ws.Range("M64").Value = -3223

$ws.Range("H65").Value = 2349.4
$ws.Range("I65").Value = 2842
$ws.Range("J65").Value = 1200
$ws.Range("K65").Value = 14210
$ws.Range("L65").Value = 6000
$ws.Range("M65").Value = -11090
$ws.Range("N65").Value = -12240

$ws.Range("H67").Value = 3478.25
$ws.Range("I67").Value = 3471
$ws.Range("K67").Value = 3471
$ws.Range("M67").Value = -2613

$ws.Range("H74").Value = 3765.2415
$ws.Range("I74").Value = 3665.3333
$ws.Range("J74").Value = 3791.3044
$ws.Range("K74").Value = 3665.3333
$ws.Range("L74").Value = 3791.3044
$ws.Range("M74").Value = -2729.3333
$ws.Range("N74").Value = -5663.3044

$ws.Range("H77").Value = 3765.2415
$ws.Range("I77").Value = 3665.3333
$ws.Range("J77").Value = 3791.3044
$ws.Range("K77").Value = 18326.6665
$ws.Range("L77").Value = 18956.522
$ws.Range("M77").Value = -13646.6665
$ws.Range("N77").Value = -28316.522

$ws.Range("H87").Value = 37074
$ws.Range("J87").Value = 37074
$ws.Range("L87").Value = 37074
$ws.Range("N87").Value = -39570

$ws.Range("H90").Value = 37074
$ws.Range("J90").Value = 37074
$ws.Range("L90").Value = 111222
$ws.Range("N90").Value = -123702

$ws.Range("H138").Value = 2113.2727
$ws.Range("I138").Value = 1471.9565
$ws.Range("J138").Value = 2815.6667
$ws.Range("K138").Value = 4415.8695
$ws.Range("L138").Value = 8447.000100000001
$ws.Range("M138").Value = 724.1305000000002
$ws.Range("N138").Value = -18727.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6916
$ws.Range("I32").Value = 7323.197
$ws.Range("K32").Value = 7323.197
$ws.Range("M32").Value = -7036.197

$ws.Range("H45").Value = 1386.5
$ws.Range("I45").Value = 1384.4
$ws.Range("J45").Value = 1391.2727
$ws.Range("K45").Value = 1384.4
$ws.Range("L45").Value = 1391.2727
$ws.Range("M45").Value = -1007.4
$ws.Range("N45").Value = -2145.2727

$ws.Range("H63").Value = 2736
$ws.Range("I63").Value = 2736
$ws.Range("K63").Value = 2736
$ws.Range("M63").Value = -2050

$ws.Range("H66").Value = 2736
$ws.Range("I66").Value = 2736
$ws.Range("K66").Value = 13680
$ws.Range("M66").Value = -10248

$ws.Range("H102").Value = 2000
$ws.Range("I102").Value = 2000
$ws.Range("K102").Value = 2000
$ws.Range("M102").Value = -378

$ws.Range("H110").Value = 600
$ws.Range("I110").Value = 600
$ws.Range("K110").Value = 600
$ws.Range("M110").Value = 1445

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 14707584
$ws.Range("I86").Value = 1684.7727
$ws.Range("J86").Value = 41668400
$ws.Range("K86").Value = 1684.7727
$ws.Range("L86").Value = 41668400
$ws.Range("M86").Value = -561.7727
$ws.Range("N86").Value = -41670646

$ws.Range("H89").Value = 14707584
$ws.Range("I89").Value = 1684.7727
$ws.Range("J89").Value = 41668400
$ws.Range("K89").Value = 8423.863499999999
$ws.Range("L89").Value = 208342000
$ws.Range("M89").Value = -2807.863499999999
$ws.Range("N89").Value = -208353232

$ws.Range("H94").Value = 1583.5
$ws.Range("I94").Value = 1445.4286
$ws.Range("J94").Value = 1690.8889
$ws.Range("K94").Value = 1445.4286
$ws.Range("L94").Value = 1690.8889
$ws.Range("M94").Value = -994.4286
$ws.Range("N94").Value = -2592.8889

$ws.Range("H99").Value = 1041.3334
$ws.Range("I99").Value = 1024
$ws.Range("K99").Value = 1024
$ws.Range("M99").Value = 474

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 29507778
$ws.Range("I4").Value = 33342668
$ws.Range("K4").Value = 33342668
$ws.Range("M4").Value = -33342556

$ws.Range("H62").Value = 2282.2222
$ws.Range("I62").Value = 2192.5
$ws.Range("J62").Value = 3000
$ws.Range("K62").Value = 2192.5
$ws.Range("L62").Value = 3000
$ws.Range("M62").Value = -1568.5
$ws.Range("N62").Value = -4248

$ws.Range("H65").Value = 2282.2222
$ws.Range("I65").Value = 2192.5
$ws.Range("J65").Value = 3000
$ws.Range("K65").Value = 10962.5
$ws.Range("L65").Value = 15000
$ws.Range("M65").Value = -7842.5
$ws.Range("N65").Value = -21240

$ws.Range("H140").Value = 30753.3
$ws.Range("J140").Value = 30753.3
$ws.Range("L140").Value = 30753.3
$ws.Range("N140").Value = -41113.3

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("N43").ClearContents()
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0

$ws.Range("H107").Value = 910.3889
$ws.Range("I107").Value = 217.16667
$ws.Range("J107").Value = 1257
$ws.Range("K107").Value = 651.50001
$ws.Range("L107").Value = 3771
$ws.Range("M107").Value = 1268.49999
$ws.Range("N107").Value = -7611

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 25643788
$ws.Range("I80").Value = 37039650
$ws.Range("K80").Value = 37039650
$ws.Range("M80").Value = -37038652

$ws.Range("H83").Value = 25643788
$ws.Range("I83").Value = 37039650
$ws.Range("K83").Value = 185198250
$ws.Range("M83").Value = -185193258

$ws.Range("H113").Value = 91883.82000000001
$ws.Range("I113").Value = 250993
$ws.Range("J113").Value = 964.2857
$ws.Range("K113").Value = 250993
$ws.Range("L113").Value = 964.2857
$ws.Range("M113").Value = -248823
$ws.Range("N113").Value = -5304.2857

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 6667483.5

$ws.Range("H16").Value = 5380
$ws.Range("I16").Value = 5333.3335
$ws.Range("J16").Value = 5450
$ws.Range("K16").Value = 5333.3335
$ws.Range("L16").Value = 5450
$ws.Range("M16").Value = -5163.3335
$ws.Range("N16").Value = -5790

$ws.Range("H22").Value = 1014.7273
$ws.Range("I22").Value = 393.33334
$ws.Range("J22").Value = 1444.9231
$ws.Range("K22").Value = 393.33334
$ws.Range("L22").Value = 1444.9231
$ws.Range("M22").Value = -98.33334000000002
$ws.Range("N22").Value = -2034.9231

$ws.Range("H27").Value = 1014.7273
$ws.Range("I27").Value = 393.33334
$ws.Range("J27").Value = 1444.9231
$ws.Range("K27").Value = 393.33334
$ws.Range("L27").Value = 1444.9231
$ws.Range("M27").Value = -286.33334
$ws.Range("N27").Value = -1658.9231

$ws.Range("H68").Value = 1821
$ws.Range("I68").Value = 2475
$ws.Range("J68").Value = 1385
$ws.Range("K68").Value = 2475
$ws.Range("L68").Value = 1385
$ws.Range("M68").Value = -1726
$ws.Range("N68").Value = -2883

$ws.Range("H71").Value = 1821
$ws.Range("I71").Value = 2475
$ws.Range("J71").Value = 1385
$ws.Range("K71").Value = 12375
$ws.Range("L71").Value = 6925
$ws.Range("M71").Value = -8631
$ws.Range("N71").Value = -14413

$ws.Range("H93").Value = 1850
$ws.Range("I93").Value = 1725
$ws.Range("J93").Value = 2000
$ws.Range("K93").Value = 1725
$ws.Range("L93").Value = 2000
$ws.Range("M93").Value = -477
$ws.Range("N93").Value = -4496

$ws.Range("H132").Value = 12828672
$ws.Range("I132").Value = 6546
$ws.Range("J132").Value = 19239736
$ws.Range("K132").Value = 19638
$ws.Range("L132").Value = 57719208
$ws.Range("M132").Value = -17108
$ws.Range("N132").Value = -57724268

$ws.Range("H136").Value = 17863628
$ws.Range("I136").Value = 33335800
$ws.Range("J136").Value = 11120
$ws.Range("K136").Value = 100007400
$ws.Range("L136").Value = 33360
$ws.Range("M136").Value = -100004850
$ws.Range("N136").Value = -38460

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("M4").ClearContents()
$ws.Range("H4").Value = 572.3333
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 572.3333
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 572.3333

$ws.Range("H122").Value = 2808.65
$ws.Range("I122").Value = 2948.6667
$ws.Range("K122").Value = 8846.000100000001
$ws.Range("M122").Value = -6396.000100000001

$ws.Range("H132").Value = 2270.1155
$ws.Range("I132").Value = 1435.8948
$ws.Range("J132").Value = 4534.4287
$ws.Range("K132").Value = 4307.6844
$ws.Range("L132").Value = 13603.2861
$ws.Range("M132").Value = -1777.6844
$ws.Range("N132").Value = -18663.2861

$ws.Range("H136").Value = 1632.1111
$ws.Range("I136").Value = 1734.9166
$ws.Range("J136").Value = 1426.5
$ws.Range("K136").Value = 5204.7498
$ws.Range("L136").Value = 4279.5
$ws.Range("M136").Value = -2654.7498
$ws.Range("N136").Value = -9379.5
